$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 1860
$ws.Range("I3").Value = 1977
$ws.Range("I4").Value = 497
$ws.Range("I5").Value = 173
$ws.Range("I6").Value = 2397
$ws.Range("I7").Value = 6904

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 70
$ws.Range("I4").Value = 32
$ws.Range("I6").Value = 43
$ws.Range("I7").Value = 240
$ws.Range("I8").Value = 426
$ws.Range("I10").Value = 54
$ws.Range("I11").Value = 121
$ws.Range("I13").Value = 9
$ws.Range("I15").Value = 86
$ws.Range("I18").Value = 56
$ws.Range("I19").Value = 202
$ws.Range("I21").Value = 48
$ws.Range("I27").Value = 63
$ws.Range("I29").Value = 447
$ws.Range("I30").Value = 22
$ws.Range("I33").Value = 325
$ws.Range("I34").Value = 31
$ws.Range("I35").Value = 9
$ws.Range("I36").Value = 89
$ws.Range("I37").Value = 222
$ws.Range("I41").Value = 30
$ws.Range("I42").Value = 228
$ws.Range("I44").Value = 54
$ws.Range("I47").Value = 53
$ws.Range("I48").Value = 69
$ws.Range("I52").Value = 144
$ws.Range("I54").Value = 152
$ws.Range("I60").Value = 38
$ws.Range("I63").Value = 26
$ws.Range("I64").Value = 71
$ws.Range("I65").Value = 161
$ws.Range("I67").Value = 266
$ws.Range("I72").Value = 25
$ws.Range("I73").Value = 66
$ws.Range("I75").Value = 26
$ws.Range("I76").Value = 112
$ws.Range("I79").Value = 174
$ws.Range("I83").Value = 130
$ws.Range("I85").Value = 328
$ws.Range("I89").Value = 69
$ws.Range("I91").Value = 78
$ws.Range("I92").Value = 20
$ws.Range("I96").Value = 92
$ws.Range("I98").Value = 45
$ws.Range("I99").Value = 123
$ws.Range("I100").Value = 8
$ws.Range("I101").Value = 6904

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 79
$ws.Range("I3").Value = 127
$ws.Range("I5").Value = 11
$ws.Range("I6").Value = 94
$ws.Range("I7").Value = 328

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 59
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 26
$ws.Range("I7").Value = 144

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 52
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 127
$ws.Range("I3").Value = 117
$ws.Range("I7").Value = 426

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 81
$ws.Range("I7").Value = 240

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 17
$ws.Range("I3").Value = 15
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 92

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 22

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I3").Value = 68
$ws.Range("I4").Value = 17
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 92
$ws.Range("I7").Value = 266

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 42
$ws.Range("I3").Value = 46
$ws.Range("I7").Value = 161

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 45
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 79
$ws.Range("I3").Value = 111
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 325

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 31
$ws.Range("I4").Value = 11
$ws.Range("I6").Value = 77
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 134
$ws.Range("I3").Value = 144
$ws.Range("I5").Value = 13
$ws.Range("I6").Value = 142
$ws.Range("I7").Value = 447

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 84
$ws.Range("I3").Value = 50
$ws.Range("I7").Value = 202

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I2").Value = 17
$ws.Range("I7").Value = 54

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I3").Value = 17
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 112

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 43

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I2").Value = 10
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 79
$ws.Range("I7").Value = 228

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("I4").Value = 3
$ws.Range("I6").Value = 9

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 54

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I2").Value = 26
$ws.Range("I7").Value = 78

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I3").Value = 4
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 48

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 49
$ws.Range("I3").Value = 49
$ws.Range("I4").Value = 9
$ws.Range("I7").Value = 174

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 14
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 71

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 28
$ws.Range("I7").Value = 89

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("I3").Value = 4
$ws.Range("I6").Value = 8

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("I4").Value = 2
$ws.Range("I6").Value = 11
$ws.Range("I7").Value = 31

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I3").Value = 16
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 18
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 9

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 66

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I3").Value = 25
$ws.Range("I4").Value = 9
$ws.Range("I7").Value = 70

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I6").Value = 11
$ws.Range("I7").Value = 20

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I3").Value = 11
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I2").Value = 10
$ws.Range("I7").Value = 26

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 25

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I6").Value = 13
$ws.Range("I7").Value = 32
